$d = $word.ActiveDocument

# 1. Christmas tone: "Tone 8a" -> "Tone 8"
$d.Content.Find.Execute("Tone 8a", $true, $false, $false, $false, $false, $true, 1, $false, "Tone 8", 2)

# 2. Epiphany tone: ", unnamed A" -> ", 11"
$d.Content.Find.Execute(", unnamed A", $true, $false, $false, $false, $false, $true, 1, $false, ", 11", 2)

# 3. Easter tone: ", 8a" -> ", 8"
$d.Content.Find.Execute(", 8a", $true, $false, $false, $false, $false, $true, 1, $false, ", 8", 2)

# 4. Ordinary Time tones part 1: ", 6a, 6b, 8b" -> ", 6, 9"
$d.Content.Find.Execute(", 6a, 6b, 8b", $true, $false, $false, $false, $false, $true, 1, $false, ", 6, 9", 2)

# 5. Ordinary Time tones part 2: ", unnamed B" -> ", 10, 12"
$d.Content.Find.Execute(", unnamed B", $true, $false, $false, $false, $false, $true, 1, $false, ", 10, 12", 2)
